# Update Betfair Back/Lay odds data on Sheet1
# This script applies targeted cell value updates matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 2.24
$ws.Range("I2").Value = 3.15
$ws.Range("K2").Value = 4.6
$ws.Range("L2").Value = 1.23
$ws.Range("Q2").Value = 1.43
$ws.Range("R2").Value = 1.8
$ws.Range("S2").Value = 2.08
$ws.Range("T2").Value = 1.44
$ws.Range("V2").Value = 1.46
$ws.Range("AI2").Value = 1000

# Row 3
$ws.Range("F3").Value = 1.98
$ws.Range("G3").Value = 2.08
$ws.Range("I3").Value = 4
$ws.Range("L3").Value = 1.27
$ws.Range("N3").Value = 4.6
$ws.Range("O3").Value = 1.24
$ws.Range("P3").Value = 2.28
$ws.Range("Q3").Value = 1.71
$ws.Range("R3").Value = 1.51
$ws.Range("S3").Value = 2.78
$ws.Range("U3").Value = 2.34
$ws.Range("W3").Value = 1.92
$ws.Range("X3").Value = 23
$ws.Range("AJ3").Value = 25
$ws.Range("AN3").Value = 12
$ws.Range("AO3").Value = 980

# Row 4
$ws.Range("F4").Value = 1.99
$ws.Range("G4").Value = 2.24
$ws.Range("Q4").Value = 2.18

# Row 5
$ws.Range("G5").Value = 3.05
$ws.Range("H5").Value = 2.48
$ws.Range("J5").Value = 3.45
$ws.Range("P5").Value = 2.06

# Row 6
$ws.Range("F6").Value = 10.5
$ws.Range("G6").Value = 27
$ws.Range("H6").Value = 1.15
$ws.Range("I6").Value = 1.21
$ws.Range("K6").Value = 15
$ws.Range("N6").Value = 6.6
$ws.Range("O6").Value = 1.12
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 1.3
$ws.Range("R6").Value = 1.92
$ws.Range("S6").Value = 1.76
$ws.Range("T6").Value = 2.1
$ws.Range("U6").Value = 1.73
$ws.Range("Z6").Value = 11.5
$ws.Range("AA6").Value = 11.5
$ws.Range("AB6").Value = 90
$ws.Range("AO6").Value = 3.3

# Row 7
$ws.Range("F7").Value = 2.24
$ws.Range("G7").Value = 2.28
$ws.Range("Z7").Value = 27

# Row 8
$ws.Range("F8").Value = 3.75
$ws.Range("H8").Value = 2.18
$ws.Range("I8").Value = 2.32
$ws.Range("P8").Value = 1.67

# Row 9
$ws.Range("G9").Value = 2.44
$ws.Range("S9").Value = 2.96
$ws.Range("T9").Value = 1.66
$ws.Range("X9").Value = 18

# Row 10
$ws.Range("H10").Value = 3.75
$ws.Range("K10").Value = 3.15

# Row 11
$ws.Range("T11").Value = 1.9
$ws.Range("Y11").Value = 10.5

# Row 12
$ws.Range("J12").Value = 3.25
$ws.Range("K12").Value = 3.85
